$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing row values ---
$ws.Range("A2").Value = "ishmuli"
$ws.Range("B2").Value = 530
$ws.Range("C2").Value = 363
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 407
$ws.Range("F2").Value = 520
$ws.Range("G2").Value = 250
$ws.Range("H2").Value = 1325
$ws.Range("I2").Value = 76.22130102040816
$ws.Range("J2").Value = 17495
$ws.Range("K2").Value = 795
$ws.Range("L2").Value = 170
$ws.Range("M2").Value = 78
$ws.Range("N2").Value = 89
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 5510
$ws.Range("Q2").Value = 525
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 530
$ws.Range("U2").Value = 0.009435015273420415
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = -321.6560824573816
$ws.Range("X2").Value = 1.002206953444424
$ws.Range("Y2").Value = 1
$ws.Range("Z2").Value = 0.004413906888848709
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 780
$ws.Range("AC2").Value = 19477
$ws.Range("AD2").Value = 6697
$ws.Range("AE2").Value = 1325
$ws.Range("AF2").Value = 2
$ws.Range("AG2").Value = 0.9565336249316566

# --- Row 3: new row ---
$ws.Range("A3").Value = "yaniv33martin"
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = 11
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 12
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 6
$ws.Range("K3").Value = 9
$ws.Range("L3").Value = 11
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 457
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 11
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = 0.9975212478233336
$ws.Range("Z3").Value = 0
$ws.Range("AA3").Value = 1
$ws.Range("AB3").Value = 11
$ws.Range("AC3").Value = 38
$ws.Range("AD3").Value = 480
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0.4
